$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -831
$ws.Range("N13").Value = -3338

$ws.Range("H62").Value = 55558556
$ws.Range("J62").Value = 6000
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 55558556
$ws.Range("J65").Value = 6000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240

$ws.Range("H98").Value = 6548.1
$ws.Range("I98").Value = 8655
$ws.Range("J98").Value = 1632
$ws.Range("K98").Value = 8655
$ws.Range("L98").Value = 1632
$ws.Range("M98").Value = -7157
$ws.Range("N98").Value = -4628

$ws.Range("H122").Value = 6548.1
$ws.Range("I122").Value = 8655
$ws.Range("J122").Value = 1632
$ws.Range("K122").Value = 25965
$ws.Range("L122").Value = 4896
$ws.Range("M122").Value = -23515
$ws.Range("N122").Value = -9796

$ws.Range("H134").Value = 35005.555
$ws.Range("J134").Value = 35005.555
$ws.Range("L134").Value = 35005.555
$ws.Range("N134").Value = -45145.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2146.8667
$ws.Range("I45").Value = 2440.3
$ws.Range("J45").Value = 1560
$ws.Range("K45").Value = 2440.3
$ws.Range("L45").Value = 1560
$ws.Range("M45").Value = -2063.3
$ws.Range("N45").Value = -2314

$ws.Range("H122").Value = 2909.2666
$ws.Range("I122").Value = 2619.182
$ws.Range("J122").Value = 3707
$ws.Range("K122").Value = 7857.545999999999
$ws.Range("L122").Value = 11121
$ws.Range("M122").Value = -5407.545999999999
$ws.Range("N122").Value = -16021

$ws.Range("H132").Value = 2870.6155
$ws.Range("I132").Value = 2116.25
$ws.Range("K132").Value = 6348.75
$ws.Range("M132").Value = -3818.75

$ws.Range("H133").Value = 34240
$ws.Range("J133").Value = 34240
$ws.Range("L133").Value = 34240
$ws.Range("N133").Value = -39300

$ws.Range("H134").Value = 32633.334
$ws.Range("J134").Value = 32633.334
$ws.Range("L134").Value = 32633.334
$ws.Range("N134").Value = -42773.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13079.667
$ws.Range("I134").Value = 34240
$ws.Range("J134").Value = 2499.5
$ws.Range("K134").Value = 102720
$ws.Range("L134").Value = 7498.5
$ws.Range("M134").Value = -100185
$ws.Range("N134").Value = -12568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 111112740
$ws.Range("J16").Value = 1262.5
$ws.Range("K16").Value = 111112740
$ws.Range("L16").Value = 1262.5
$ws.Range("M16").Value = -111112453
$ws.Range("N16").Value = -1836.5

$ws.Range("H31").Value = 1405.28
$ws.Range("I31").Value = 1392.1
$ws.Range("J31").Value = 1458
$ws.Range("K31").Value = 1392.1
$ws.Range("L31").Value = 1458
$ws.Range("M31").Value = -1097.1
$ws.Range("N31").Value = -2048

$ws.Range("H34").Value = 1405.28
$ws.Range("I34").Value = 1392.1
$ws.Range("J34").Value = 1458
$ws.Range("K34").Value = 1392.1
$ws.Range("L34").Value = 1458
$ws.Range("M34").Value = -1190.1
$ws.Range("N34").Value = -1862

$ws.Range("H105").Value = 1097.5
$ws.Range("I105").Value = 963.3333
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 963.3333
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 783.6667
$ws.Range("N105").Value = -4994

$ws.Range("I113").Value = 111112740
$ws.Range("J113").Value = 1262.5
$ws.Range("K113").Value = 111112740
$ws.Range("L113").Value = 1262.5
$ws.Range("M113").Value = -111110570
$ws.Range("N113").Value = -5602.5

$ws.Range("H132").Value = 2783.4546
$ws.Range("I132").Value = 1960
$ws.Range("J132").Value = 3254
$ws.Range("K132").Value = 5880
$ws.Range("L132").Value = 9762
$ws.Range("M132").Value = -3350
$ws.Range("N132").Value = -14822

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 500
$ws.Range("I57").Value = 500
$ws.Range("K57").Value = 1500
$ws.Range("M57").Value = -941

$ws.Range("H86").Value = 613.2778
$ws.Range("I86").Value = 638.8889
$ws.Range("J86").Value = 587.6667
$ws.Range("K86").Value = 1916.6667
$ws.Range("L86").Value = 1763.0001
$ws.Range("M86").Value = -730.6667000000002
$ws.Range("N86").Value = -4135.0001

$ws.Range("H89").Value = 613.2778
$ws.Range("I89").Value = 638.8889
$ws.Range("J89").Value = 587.6667
$ws.Range("K89").Value = 5750.0001
$ws.Range("L89").Value = 5289.0003
$ws.Range("M89").Value = 177.9998999999998
$ws.Range("N89").Value = -17145.0003

$ws.Range("H107").Value = 3694.9333
$ws.Range("J107").Value = 4511.2085
$ws.Range("L107").Value = 13533.6255
$ws.Range("N107").Value = -17373.6255

$ws.Range("H117").Value = 947.13336
$ws.Range("I117").Value = 575.5
$ws.Range("J117").Value = 1690.4
$ws.Range("K117").Value = 1726.5
$ws.Range("L117").Value = 5071.200000000001
$ws.Range("M117").Value = 1715.5
$ws.Range("N117").Value = -11955.2

$ws.Range("H120").Value = 5825.8
$ws.Range("I120").Value = 5000
$ws.Range("J120").Value = 6032.25
$ws.Range("K120").Value = 15000
$ws.Range("L120").Value = 18096.75
$ws.Range("M120").Value = -10162
$ws.Range("N120").Value = -27772.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4653.222
$ws.Range("I122").Value = 4734.875
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 14204.625
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -11754.625
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1466.6666
$ws.Range("I22").Value = 1200
$ws.Range("K22").Value = 1200
$ws.Range("M22").Value = -905

$ws.Range("H27").Value = 1466.6666
$ws.Range("I27").Value = 1200
$ws.Range("K27").Value = 1200
$ws.Range("M27").Value = -1093

$ws.Range("H82").Value = 2021.5625
$ws.Range("I82").Value = 2014.7273
$ws.Range("K82").Value = 2014.7273
$ws.Range("M82").Value = -1653.7273

$ws.Range("H85").Value = 2021.5625
$ws.Range("I85").Value = 2014.7273
$ws.Range("K85").Value = 2014.7273
$ws.Range("M85").Value = -766.7273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11906086
$ws.Range("I122").Value = 15626350
$ws.Range("K122").Value = 46879050
$ws.Range("M122").Value = -46876600
